$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update the Moisturizer data row first: Aloe -> ALOE, Almond -> ALMOND
$ws.Range("C7").Value = "ALOE"
$ws.Range("D7").Value = "ALMOND"

# Copy style from C6/D6 into new E6/F6 header cells, then set their text
$ws.Range("C6").Copy($ws.Range("E6"))
$ws.Range("D6").Copy($ws.Range("F6"))
$ws.Range("E6").Value = "ProductContent3"
$ws.Range("F6").Value = "ProductContent4"

# add SPF-50/SPF-30 to the Moisturizer row
$ws.Range("E7").Value = "SPF-50"
$ws.Range("F7").Value = "SPF-30"

# Update selection to F6
$ws.Range("F6").Select() | Out-Null
